$d = $word.ActiveDocument

# 1) Insert a bare empty paragraph between the "Hudson Valley ... 2015 - 2017"
#    paragraph and the "Skills" Heading1 paragraph.
$eduPara = $d.Paragraphs.Item(6)
$insertPoint = $d.Range($eduPara.Range.End, $eduPara.Range.End)
$insertPoint.InsertXML("<w:p/>")

# 2) Remove the stray _GoBack bookmark that sits between "simulate " and
#    " different levels of traffic".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3) Merge the three runs that spell out "|  Survival Horror  ( C++ / Unreal
#    Engine ) |  " into a single run with identical text/formatting.
$d.Content.Find.Execute("|  Survival Horror  ( C++ / Unreal Engine ) |  ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "|  Survival Horror  ( C++ / Unreal Engine ) |  ", 2)

# 4) Add a new empty ListBullet-styled paragraph (no numbering, hanging
#    indent) right after "Implemented dynamically-loaded dialogue trees..."
#    that carries the _GoBack bookmark, before the "Activities" heading.
$dlgPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Implemented dynamically-loaded dialogue trees*") {
        $dlgPara = $p
        break
    }
}
$activitiesPara = $dlgPara.Next()
$activitiesStart = $activitiesPara.Range.Start
$insertPoint2 = $d.Range($activitiesStart, $activitiesStart)
# Insert the new paragraph followed by a throwaway empty one -- a lone
# <w:p> merges into the following (Activities) paragraph, but a second
# trailing <w:p/> forces a clean split, leaving Activities untouched.
$newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='0'/></w:numPr></w:pPr>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p><w:p/>"
$insertPoint2.InsertXML($newParaXml)

$newPara = $dlgPara.Next()
$newPara.LeftIndent = 7.2
$newPara.FirstLineIndent = -7.2

# Drop the throwaway paragraph that InsertXML left behind right before
# "Activities".
$throwaway = $newPara.Next()
$throwaway.Range.Delete()
